$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The attendance app re-uploaded the log: the first excuse row now reflects a
# different student id and the following day's date, and the rest of the
# previously-logged rows (which belonged to the prior upload) are gone.

# Student ID for row 2 changed 200933 -> 200850.
# Write it through a text formula + paste-as-values so it lands back in the
# sheet as a literal text value (matching the original inlineStr cell type)
# instead of being auto-coerced to a number.
$ws.Range("A2").Formula = "=""200850"""
$ws.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues

# Log Date for row 2 changed 25/10/2025 -> 26/10/2025 (also kept as text).
$ws.Range("C2").Formula = "=""26/10/2025"""
$ws.Range("C2").Copy()
$ws.Range("C2").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false

# Rows 3-7 (the other excuse entries from the previous upload) are removed,
# shrinking the sheet down to just the header + the single remaining entry.
$ws.Rows("3:7").Delete()
